# "Generate Report for Archive"
#
# 1. Status text "Ready for handoff" -> "In Translation" (Overview sheet
#    columns E/F, and the Status column (C) on each per-locale sheet).
# 2. Narrow the "Status" column(s) from ~17.22 chars down to ~13.41 chars
#    (Overview!E:F and the "Status" column on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

# --- 1. Replace the status text everywhere it appears on every sheet ---
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    [void]$ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Resize the "Status" columns ---
# Target stored column width ~= 13.4101845877511 characters. The
# ColumnWidth COM setter here only lands on an MDW=6 pixel grid
# (stored = (Round(ColumnWidth*6)+5)/6), so 12.5 is the closest input
# that reproduces the narrower width seen in the target workbook.
$newStatusWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1").EntireColumn.ColumnWidth = $newStatusWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $newStatusWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $newStatusWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $newStatusWidth
